# Update bus voltage magnitude results (vm_pu) for the 380 kV case (Case_2_251).
# New per-unit voltage values for rows 2-25 (bus results), columns B-F and I-N.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> column letter -> new value
$vmPuUpdates = @{
    2 = @{ "B"=1.02; "C"=1.026050530594341; "D"=1.036784790796466; "E"=1.046692544700521; "F"=1.050358477651043; "I"=1.032065972700086; "J"=1.031216133384383; "K"=1.039577541219736; "L"=1.049457268143794; "M"=1.05311297624636; "N"=1.014429922827338 }
    3 = @{ "B"=1.02; "C"=1.026900658436384; "D"=1.037452186196245; "E"=1.047592242606375; "F"=1.05123218634374; "I"=1.032178456947; "J"=1.031706368123838; "K"=1.040054986912077; "L"=1.050168431362989; "M"=1.053798957428911; "N"=1.014593736074927 }
    4 = @{ "B"=1.02; "C"=1.0274513327783; "D"=1.037884373623767; "E"=1.048175716322523; "F"=1.051798546887034; "I"=1.032249987199211; "J"=1.032023526314415; "K"=1.040363582255681; "L"=1.050629274882451; "M"=1.054243204433859; "N"=1.014699671813164 }
    5 = @{ "B"=1.02; "C"=1.027682974772709; "D"=1.038066144356828; "E"=1.048421320262802; "F"=1.052036885427652; "I"=1.032279757514451; "J"=1.032156844904049; "K"=1.040493232056643; "L"=1.050823173178656; "M"=1.054430053189208; "N"=1.014744191817669 }
    6 = @{ "B"=1.02; "C"=1.027721876561164; "D"=1.038096669051852; "E"=1.04846257646487; "F"=1.052076917570649; "I"=1.032284738411414; "J"=1.032179228784311; "K"=1.040514995882612; "L"=1.050855738875253; "M"=1.054461430991709; "N"=1.014751666007517 }
    7 = @{ "B"=1.02; "C"=1.027454427447622; "D"=1.037886802144471; "E"=1.048178996872303; "F"=1.051801730635332; "I"=1.032250386175173; "J"=1.03202530778325; "K"=1.040365314973162; "L"=1.050631865135452; "M"=1.054245700774968; "N"=1.014700266752633 }
    8 = @{ "B"=1.02; "C"=1.026337713295068; "D"=1.037010269150452; "E"=1.046996329704573; "F"=1.050653540817656; "I"=1.032104246652791; "J"=1.031381821482974; "K"=1.039738966637268; "L"=1.049697468742765; "M"=1.053344728967829; "N"=1.014485296931905 }
    9 = @{ "B"=1.02; "C"=1.024374465252139; "D"=1.035468372878139; "E"=1.0449224236975; "F"=1.048638120734441; "I"=1.031837150527536; "J"=1.030247536773591; "K"=1.038632685867692; "L"=1.04805616748194; "M"=1.051760015100542; "N"=1.014106034459082 }
    10 = @{ "B"=1.02; "C"=1.023068777906863; "D"=1.034442340707651; "E"=1.043546720202966; "F"=1.047299886760564; "I"=1.031652684636877; "J"=1.029491162818276; "K"=1.037893513862206; "L"=1.046965570168699; "M"=1.050705592108982; "N"=1.013852911172504 }
    11 = @{ "B"=1.02; "C"=1.022504165361243; "D"=1.033998530072781; "E"=1.042952684607432; "F"=1.046721715825548; "I"=1.031571298637296; "J"=1.029163615812932; "K"=1.037573069303482; "L"=1.046494203335785; "M"=1.050249522035635; "N"=1.013743244751682 }
    12 = @{ "B"=1.02; "C"=1.022294558589736; "D"=1.033833751161605; "E"=1.042732283516555; "F"=1.046507153561802; "I"=1.031540841875713; "J"=1.029041946644877; "K"=1.037453986448963; "L"=1.046319248806028; "M"=1.050080194321471; "N"=1.013702500834995 }
    13 = @{ "B"=1.02; "C"=1.022339514704772; "D"=1.033869093505018; "E"=1.042779548948371; "F"=1.046553168988754; "I"=1.031547385188228; "J"=1.029068045239345; "K"=1.037479532614413; "L"=1.046356771144002; "M"=1.050116512180184; "N"=1.013711240940173 }
    14 = @{ "B"=1.02; "C"=1.022486836836038; "D"=1.03398490791498; "E"=1.042934461070407; "F"=1.046703976029351; "I"=1.031568785686208; "J"=1.029153558668608; "K"=1.03756322700154; "L"=1.046479738836064; "M"=1.050235523773617; "N"=1.013739877026316 }
    15 = @{ "B"=1.02; "C"=1.022577622208302; "D"=1.034056274605554; "E"=1.043029940747103; "F"=1.046796919274125; "I"=1.03158194126127; "J"=1.029206245869914; "K"=1.037614786564556; "L"=1.046555520840212; "M"=1.05030886099741; "N"=1.013757519494738 }
    16 = @{ "B"=1.02; "C"=1.023106265493226; "D"=1.03447180497294; "E"=1.043586179373167; "F"=1.047338285467764; "I"=1.031658054176187; "J"=1.029512900460153; "K"=1.037914772858566; "L"=1.04699687167344; "M"=1.050735870690983; "N"=1.013860188095652 }
    17 = @{ "B"=1.02; "C"=1.023438073627928; "D"=1.034732582778479; "E"=1.043935537032603; "F"=1.047678217729828; "I"=1.0317053936975; "J"=1.029705249034772; "K"=1.038102846027929; "L"=1.047273953135197; "M"=1.051003858124187; "N"=1.013924573013205 }
    18 = @{ "B"=1.02; "C"=1.023631684802228; "D"=1.034884735101551; "E"=1.044139470822353; "F"=1.047876619117005; "I"=1.031732860201767; "J"=1.029817439527553; "K"=1.038212509344554; "L"=1.047435653806465; "M"=1.051160218992789; "N"=1.013962121598945 }
    19 = @{ "B"=1.02; "C"=1.023697713563801; "D"=1.034936622696781; "E"=1.04420903395241; "F"=1.047944289948641; "I"=1.031742200813856; "J"=1.029855693011182; "K"=1.038249895455914; "L"=1.047490803713926; "M"=1.051213542157924; "N"=1.013974923657061 }
    20 = @{ "B"=1.02; "C"=1.023402466183172; "D"=1.034704599116129; "E"=1.043898037763721; "F"=1.047641733308537; "I"=1.031700329694922; "J"=1.029684612173014; "K"=1.038082671330822; "L"=1.047244216243616; "M"=1.050975100598032; "N"=1.013917665743179 }
    21 = @{ "B"=1.02; "C"=1.022443450940901; "D"=1.033950801452573; "E"=1.042888836390772; "F"=1.046659561684499; "I"=1.031562490016927; "J"=1.02912837718556; "K"=1.037538582611228; "L"=1.046443524256656; "M"=1.050200475686788; "N"=1.013731444655488 }
    22 = @{ "B"=1.02; "C"=1.021841148789649; "D"=1.033477277091565; "E"=1.042255759791312; "F"=1.046043167278541; "I"=1.031474515307228; "J"=1.028778629783475; "K"=1.037196172126598; "L"=1.04594086247674; "M"=1.049713884648793; "N"=1.01361430864868 }
    23 = @{ "B"=1.02; "C"=1.022160376665953; "D"=1.033728261079399; "E"=1.042591227779217; "F"=1.046369821321522; "I"=1.031521276270325; "J"=1.028964039032247; "K"=1.037377720316617; "L"=1.046207259963685; "M"=1.049971792920496; "N"=1.013676409387535 }
    24 = @{ "B"=1.02; "C"=1.023418555429041; "D"=1.034717243589655; "E"=1.043914981576963; "F"=1.047658218663145; "I"=1.031702618350397; "J"=1.029693937092494; "K"=1.038091787521223; "L"=1.047257652806347; "M"=1.050988094737098; "N"=1.013920786859923 }
    25 = @{ "B"=1.02; "C"=1.024881464179338; "D"=1.03586666309527; "E"=1.045457369668054; "F"=1.049158215048448; "I"=1.031907332497924; "J"=1.03054081422287; "K"=1.038918982890455; "L"=1.048479854707184; "M"=1.052169346664247; "N"=1.01420413447273 }
}

foreach ($rowNum in $vmPuUpdates.Keys) {
    $rowValues = $vmPuUpdates[$rowNum]
    foreach ($colLetter in $rowValues.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $rowValues[$colLetter]
    }
}
